$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Sheet2" after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Header row
$ws2.Range("A1").Value = "Lights"
$ws2.Range("B1").Value = "Length"
$ws2.Range("C1").Value = "Amps"

# Data rows (entered in the order the author appears to have typed them,
# to reproduce the shared-string insertion order)
$ws2.Range("A3").Value = "Ring Light"
$ws2.Range("A2").Value = "8x8 Matric"
$ws2.Range("A5").Value = "32x8 Matrix"
$ws2.Range("A4").Value = "Strip"

$ws2.Range("B2").Value = 64
$ws2.Range("B3").Value = 16
$ws2.Range("B4").Value = 150
$ws2.Range("B5").Value = 256

$ws2.Range("C2").Formula = "=B2*0.06"
$ws2.Range("C3").Formula = "=B3*0.06"
$ws2.Range("C4").Formula = "=B4*0.06"
$ws2.Range("C5").Formula = "=B5*0.06"

# Select the full data range on Sheet1 (unchanged data, but selection moved)
$ws1.Range("A1:D11").Select() | Out-Null

# Select the full data range on the new sheet, then make Sheet2 the active/visible tab
$ws2.Range("A1:C5").Select() | Out-Null
$ws2.Activate() | Out-Null

Write-Output "done"
